# Update overall_stats values for rows 2-6 (RandomForest, XGBoost,
# Logistic Regression, Voting Classifier, Stacking Classifier) to reflect
# the new model stats (3134 7/4/2025 5Trade).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 2).Value  = 1          # B: Total Trades
    $ws.Cells.Item($row, 3).Value  = 1          # C: Total Win Count
    $ws.Cells.Item($row, 5).Value  = 46         # E: Total Profit
    $ws.Cells.Item($row, 7).Value  = 15.0345    # G: Total Fee
    $ws.Cells.Item($row, 9).Value  = 3.5        # I: Average R
    $ws.Cells.Item($row, 10).Value = 30.97      # J: Max Profit
    $ws.Cells.Item($row, 11).Value = 30.97      # K: Average Profit
    $ws.Cells.Item($row, 13).Value = 0.5        # M: Total Time (hours)
    $ws.Cells.Item($row, 14).Value = 0.5        # N: Average Time (hours)
    $ws.Cells.Item($row, 15).Value = 30.97      # O: Realized Profit/Loss
    $ws.Cells.Item($row, 16).Value = 3.06       # P: R
}
